$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update leadership bio links (F2/F10 share the Lauren/Drew string; F8 is Andi)
$ws.Range("F2").Value = "[Lauren Chenarides](https://dataifa.github.io/difa-project/lauren_chenarides.html), [Drew Hanks](https://dataifa.github.io/difa-project/drew_hanks.html)"
$ws.Range("F10").Value = "[Lauren Chenarides](https://dataifa.github.io/difa-project/lauren_chenarides.html), [Drew Hanks](https://dataifa.github.io/difa-project/drew_hanks.html)"
$ws.Range("F8").Value = "[Andi Carlson](https://dataifa.github.io/difa-project/andi_carlson.html)"

# Update the active selection to F10
$ws.Range("F10").Select()
